# Se procesan de nuevo los datos con las nuevas dimensiones curadas
# Columns C (residencia), F (sexo), I (residencia-provincia-nombre) and
# J (residencia-ccaa-nombre) move from "dimension"/"dim" semantics to
# "measure"/"medida" semantics, and their mapping-file references in row 5
# are cleared since they no longer apply.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-dimension:* -> iaest-measure:* (and refArea column becomes a
# residencia-provincia-nombre measure)
$ws.Range("C2").Value = "iaest-measure:residencia"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("I2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("J2").Value = "iaest-measure:residencia-ccaa-nombre"

# Row 3: dim -> medida
$ws.Range("C3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: skos:Concept / URI-Provincia -> xsd:int
$ws.Range("C4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: clear mapping file references for the columns that are no longer
# dimensions (fully remove the cells, not just blank their value)
$ws.Range("C5").Clear()
$ws.Range("F5").Clear()
$ws.Range("J5").Clear()
